function Set-CellText($ws, $row, $col, [string]$val) {
    $cell = $ws.Cells.Item($row, $col)
    if ($val -match '^[+-]?[0-9]*\.?[0-9]+$') {
        $cell.NumberFormat = "@"
    }
    $cell.Value = $val
}

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Price / Volume updates for rows 2-45 ---
Set-CellText $ws 2 4 "26.272.22"
Set-CellText $ws 2 5 "  +1.89%  "
Set-CellText $ws 3 4 "1.647.99"
Set-CellText $ws 3 5 "  +0.61%  "
Set-CellText $ws 4 5 "  -0.11%  "
Set-CellText $ws 5 4 "216.86"
Set-CellText $ws 5 5 "  +0.55%  "
Set-CellText $ws 6 5 "  +0.51%  "
Set-CellText $ws 7 5 "  -0.11%  "
Set-CellText $ws 8 5 "  +0.43%  "
Set-CellText $ws 9 5 "  +0.57%  "
Set-CellText $ws 10 4 "19.92"
Set-CellText $ws 10 5 "  +1.51%  "
Set-CellText $ws 11 5 "  +0.11%  "
Set-CellText $ws 12 5 "  +0.90%  "
Set-CellText $ws 13 4 "1.875.14"
Set-CellText $ws 13 5 "  +0.58%  "
Set-CellText $ws 14 4 "1.636.54"
Set-CellText $ws 14 5 "  -0.33%  "
Set-CellText $ws 15 5 "  -2.79%  "
Set-CellText $ws 16 5 "  +0.44%  "
Set-CellText $ws 17 4 "63.39"
Set-CellText $ws 17 5 "  +0.31%  "
Set-CellText $ws 18 4 "26.279.99"
Set-CellText $ws 18 5 "  +1.77%  "
Set-CellText $ws 19 5 "  -0.11%  "
Set-CellText $ws 20 5 "  -0.60%  "
Set-CellText $ws 21 4 "195.64"
Set-CellText $ws 21 5 "  +1.64%  "
Set-CellText $ws 22 5 "  +1.28%  "
Set-CellText $ws 23 4 "6.32"
Set-CellText $ws 24 5 "  -2.84%  "
Set-CellText $ws 25 4 "143.50"
Set-CellText $ws 25 5 "  +1.26%  "
Set-CellText $ws 26 5 "  -0.16%  "
Set-CellText $ws 27 5 "  +0.63%  "
Set-CellText $ws 28 5 "  +0.41%  "
Set-CellText $ws 29 5 "  +1.01%  "
Set-CellText $ws 30 5 "  +0.93%  "
Set-CellText $ws 31 5 "  +2.54%  "
Set-CellText $ws 32 5 "  +0.87%  "
Set-CellText $ws 33 5 "  +0.88%  "
Set-CellText $ws 34 5 "  +2.16%  "
Set-CellText $ws 35 5 "  +1.25%  "
Set-CellText $ws 36 5 "  +1.02%  "
Set-CellText $ws 37 4 "1.138.76"
Set-CellText $ws 37 5 "  +0.25%  "
Set-CellText $ws 38 5 "  +1.68%  "
Set-CellText $ws 39 5 "  -1.68%  "
Set-CellText $ws 40 5 "  +1.42%  "
Set-CellText $ws 41 5 "  -0.15%  "
Set-CellText $ws 42 4 "100.59"
Set-CellText $ws 42 5 "  -0.09%  "
Set-CellText $ws 43 5 "  -0.82%  "
Set-CellText $ws 44 4 "0.801"
Set-CellText $ws 44 5 "  -0.20%  "
Set-CellText $ws 45 4 "1.783.97"
Set-CellText $ws 45 5 "  +0.59%  "

# --- Rows 46-51: BabyDogeCoin inserted at rank 46 (rank index unaffected),
#     Algorand dropped off the bottom, everything else shifts down one rank ---
Set-CellText $ws 46 2 "BabyDogeCoin"
Set-CellText $ws 46 3 "https://coinranking.com/coin/JY1_q2c0g+babydogecoin-babydoge"
Set-CellText $ws 46 4 "0.0₆0112"
Set-CellText $ws 46 5 "  -0.31%  "
Set-CellText $ws 47 2 "Aave"
Set-CellText $ws 47 3 "https://coinranking.com/coin/ixgUfzmLR+aave-aave"
Set-CellText $ws 47 4 "57.22"
Set-CellText $ws 47 5 "  +3.38%  "
Set-CellText $ws 48 2 "Cronos"
Set-CellText $ws 48 3 "https://coinranking.com/coin/65PHZTpmE55b+cronos-cro"
Set-CellText $ws 48 4 "0.0517"
Set-CellText $ws 48 5 "  +2.99%  "
Set-CellText $ws 49 2 "RenderToken"
Set-CellText $ws 49 3 "https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr"
Set-CellText $ws 49 4 "1.47"
Set-CellText $ws 49 5 "  +3.45%  "
Set-CellText $ws 50 2 "EnergySwap"
Set-CellText $ws 50 3 "https://coinranking.com/coin/SbWqqTui-+energyswap-ens"
Set-CellText $ws 50 4 "7.71"
Set-CellText $ws 50 5 "  +3.35%  "
Set-CellText $ws 51 2 "Mantle"
Set-CellText $ws 51 3 "https://coinranking.com/coin/BoI4ux0nd+mantle-mnt"
Set-CellText $ws 51 4 "0.418"
Set-CellText $ws 51 5 "  +0.25%  "
